$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C9").Value = "NXPsemiconductor"
$ws.Range("C29").Value = '"Smok3y Is Back"'
$ws.Range("C50").Value = '"落日幻影 哈哈哈"'
$ws.Range("C143").Value = "rip_indragon"
